$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values in rows 2-4 (new TPM-derived numbers) ---

# Row 2
$ws.Range("G2").Value = 0.7999296666666668
$ws.Range("H2").Value = 2.399789
$ws.Range("M2").Value = 0.5471986666666667
$ws.Range("N2").Value = 1.641596
$ws.Range("O2").Value = 0.02624837670119588
$ws.Range("P2").Value = 0.02624837670119588
$ws.Range("Q2").Value = 0.4377204470271112
$ws.Range("R2").Value = 3.939484023244
$ws.Range("S2").Value = 0.02624837670119588
$ws.Range("T2").Value = 0.02624837670119588

# Row 3
$ws.Range("G3").Value = 0.7999296666666668
$ws.Range("H3").Value = 2.399789
$ws.Range("O3").Value = 0.8479060840514853
$ws.Range("P3").Value = 0.8479060840514852
$ws.Range("Q3").Value = 14.139763169854
$ws.Range("R3").Value = 127.257868528686
$ws.Range("S3").Value = 0.8479060840514853
$ws.Range("T3").Value = 0.8479060840514852

# Row 4
$ws.Range("G4").Value = 0.7999296666666668
$ws.Range("H4").Value = 2.399789
$ws.Range("M4").Value = 2.604410666666667
$ws.Range("N4").Value = 7.813232000000001
$ws.Range("O4").Value = 0.1249300417336775
$ws.Range("P4").Value = 0.1249300417336775
$ws.Range("Q4").Value = 2.083345356449778
$ws.Range("R4").Value = 18.750108208048
$ws.Range("S4").Value = 0.1249300417336775
$ws.Range("T4").Value = 0.1249300417336775

# --- Add new row 5 ---
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Gdnf"
$ws.Range("C5").Value = "Gfra1"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7999296666666668
$ws.Range("H5").Value = 2.399789
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01908533333333333
$ws.Range("N5").Value = 0.057256
$ws.Range("O5").Value = 0.0009154975136414022
$ws.Range("P5").Value = 0.0009154975136414022
$ws.Range("Q5").Value = 0.01526692433155556
$ws.Range("R5").Value = 0.137402318984
$ws.Range("S5").Value = 0.0009154975136414022
$ws.Range("T5").Value = 0.0009154975136414022
